$wb = $excel.ActiveWorkbook

# --- Sheet "OFF" (row 2) ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 199
$wsOff.Range("C2").Value = 148
$wsOff.Range("D2").Value = 53
$wsOff.Range("E2").Value = 28

# --- Sheet "DEF" (row 2) ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 250
$wsDef.Range("C2").Value = 172
$wsDef.Range("D2").Value = 62
$wsDef.Range("E2").Value = 19
$wsDef.Range("F2").Value = 4
$wsDef.Range("G2").Value = 8
